# "Add files via upload" - Anu - Cash Management Files Uploaded
#
# The uploaded version of this test-data workbook no longer ships the
# sample login/URL values that used to live on the Input_Value sheet
# (A2:C2) - they've been blanked out (formatting/styles kept) and the
# URL hyperlink that used to sit on A2 has been removed along with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Drop the hyperlink that pointed A2 at the old environment URL.
$ws.Hyperlinks.Delete()

# Clear out the old sample credentials/URL (keeps each cell's style).
$ws.Range("A2:C2").ClearContents()

# Leave the A2:C2 row selected, matching the saved selection state.
$ws.Range("A2:C2").Select()
